# Generate Report for Handback
#
# The "zh-cn"/"de-de" status column (reused on the Overview sheet as well)
# flips from "in sync" to "not in sync" with en-US, and the two localized
# sheets pick up a freshly regenerated "Correspond Handback DateTime" for
# the second file (7b262a9c-ac66-45ae-927b-19e8ea82e08f...).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: not in sync with en-US"

# Overview sheet: zh-cn / de-de status cells for both rows.
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# zh-cn sheet: Status column for both rows.
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

# de-de sheet: Status column for both rows.
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# Refresh the "Correspond Handback DateTime" for the second data row
# (7b262a9c-ac66-45ae-927b-19e8ea82e08f...) on each language sheet.
$wsZhCn.Range("K3").Value = "2016-10-20 09:35:29"
$wsDeDe.Range("K3").Value = "2016-10-20 09:35:47"

# The longer status text no longer fits the old column width, so the
# status columns grow to fit it (Overview!E:F, and the "Status" column on
# each language sheet).
$newWidth = 33.4602203369141
$wsOverview.Range("E:F").ColumnWidth = $newWidth
$wsZhCn.Range("C:C").ColumnWidth = $newWidth
$wsDeDe.Range("C:C").ColumnWidth = $newWidth
